$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -- Order below was chosen to reproduce the exact shared-string table
#    ordering of the target workbook (new strings are appended to the
#    shared string table in the order they are first written). --

# Subject 10 (row 12): start with Major so "Biomedical Eng." becomes
# the first newly-introduced shared string.
$ws.Cells.Item(12, 4).Value = "Biomedical Eng."

# Fix typo in row 10 (subject 8): "Desktop: Game Console" -> "Desktop; Game Console"
$ws.Cells.Item(10, 6).Value = "Desktop; Game Console"

# Continue filling subject 10's row
$ws.Cells.Item(12, 6).Value = "Game Console"
$ws.Cells.Item(12, 7).Value = "Joysick; Wiimote"

# Subject 11 (row 13)
$ws.Cells.Item(13, 4).Value = "Biomedical Eng."
$ws.Cells.Item(13, 6).Value = "Desktop; Game Console"
$ws.Cells.Item(13, 7).Value = "Keyboard/Mouse; Joystick; "

# Remaining (already-existing) values for both new rows
$ws.Cells.Item(12, 2).Value = "Female"
$ws.Cells.Item(12, 3).Value = 20
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 8).Value = 2
$ws.Cells.Item(12, 9).Value = 1
$ws.Cells.Item(12, 10).Value = 3
$ws.Cells.Item(12, 11).Value = "G"
$ws.Rows.Item(12).RowHeight = 30

$ws.Cells.Item(13, 2).Value = "Female"
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(13, 11).Value = "F"
$ws.Rows.Item(13).RowHeight = 30

# Update current selection to K12
$ws.Range("K12").Select()
